# Trade #17 closed at 2026-02-18 00:12:22 - unknown UNKNOWN +0.000%
#
# This script:
#  1) Updates the Summary sheet roll-up metrics.
#  2) Updates the Strategy Status row for MarketMaking.
#  3) Closes trade #47 (row 48 in "All Trades", row 19 in "MarketMaking")
#     as CLOSED / early_exit.
#  4) Appends 3 new OPEN trades (#75 EMAArbitrage, #76 momentum,
#     #77 HighProbConvergence) to "All Trades" and to each of their
#     respective per-strategy sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [int]$row, [int]$col, [string]$text)
    # Force text storage so date-like strings ("2026-02-18") are not
    # silently reinterpreted as date serials. A leading apostrophe forces
    # Excel to treat the entry as literal text; re-applying the "Normal"
    # style afterwards clears the implicit "Text" number format Excel
    # stamps on the cell, so the stored cell style matches an untouched
    # cell (only the content type changes, from blank to text).
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

function Set-EmptyTextCell {
    param($ws, [int]$row, [int]$col)
    # Source rows (e.g. an unfilled "Exit Price"/"Exit Reason" on a still-
    # OPEN trade) store an explicit-but-empty text marker rather than
    # leaving the cell truly blank. A lone leading apostrophe is the COM
    # equivalent of that: an explicit zero-length text entry.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.77
$wsSummary.Range("B4").Value = 0.87
$wsSummary.Range("B5").Value = 0.39
$wsSummary.Range("B6").Value = 45
$wsSummary.Range("B7").Value = 26
$wsSummary.Range("B9").Value = 57.78

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 99.77
$wsStatus.Range("D6").Value = 16
$wsStatus.Range("E6").Value = -0.04
$wsStatus.Range("F6").Value = -0.23
$wsStatus.Range("G6").Value = 62.5

# ---------------------------------------------------------------------
# 3) Close trade #47 (MarketMaking) in "All Trades" (row 48)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G48").Value = 0.6899999999999999
$wsAll.Range("H48").Value = "CLOSED"
$wsAll.Range("I48").Value = 6.1538
$wsAll.Range("J48").Value = 0.04
$wsAll.Range("K48").Value = 99.77
$wsAll.Range("L48").Value = "early_exit"
$wsAll.Range("M48").Value = 0.12

# ... and the mirrored row in the per-strategy "MarketMaking" sheet (row 19)
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G19").Value = 0.6899999999999999
$wsMM.Range("H19").Value = "CLOSED"
$wsMM.Range("I19").Value = 6.1538
$wsMM.Range("J19").Value = 0.04
$wsMM.Range("K19").Value = 99.77
$wsMM.Range("P19").Value = "early_exit"
$wsMM.Range("Q19").Value = 0.12

# ---------------------------------------------------------------------
# 4) Append new trades
# ---------------------------------------------------------------------

# --- "All Trades" row 76: trade #75, EMAArbitrage, UP, OPEN ---
$wsAll.Range("A76").Value = 75
Set-TextCell $wsAll 76 2 "2026-02-18"
Set-TextCell $wsAll 76 3 "00:11:14"
$wsAll.Range("D76").Value = "EMAArbitrage"
$wsAll.Range("E76").Value = "UP"
$wsAll.Range("F76").Value = 0.46
Set-EmptyTextCell $wsAll 76 7               # G76 Exit Price (unset)
$wsAll.Range("H76").Value = "OPEN"
$wsAll.Range("I76").Value = 0
$wsAll.Range("J76").Value = 0
$wsAll.Range("K76").Value = 100
Set-EmptyTextCell $wsAll 76 12              # L76 Exit Reason (unset)
$wsAll.Range("M76").Value = 0
$wsAll.Range("N76").Value = 0
$wsAll.Range("O76").Value = 0
$wsAll.Range("P76").Value = 0.7169
$wsAll.Range("Q76").Value = "EMA:up, RSI:50.0, ROC:21.69% | 2/3 UP"

# --- "All Trades" row 77: trade #76, momentum, DOWN, OPEN ---
$wsAll.Range("A77").Value = 76
Set-TextCell $wsAll 77 2 "2026-02-18"
Set-TextCell $wsAll 77 3 "00:12:16"
$wsAll.Range("D77").Value = "momentum"
$wsAll.Range("E77").Value = "DOWN"
$wsAll.Range("F77").Value = 0.65
Set-EmptyTextCell $wsAll 77 7               # G77 Exit Price (unset)
$wsAll.Range("H77").Value = "OPEN"
$wsAll.Range("I77").Value = 0
$wsAll.Range("J77").Value = 0
$wsAll.Range("K77").Value = 100
Set-EmptyTextCell $wsAll 77 12              # L77 Exit Reason (unset)
$wsAll.Range("M77").Value = 0
$wsAll.Range("N77").Value = 0
$wsAll.Range("O77").Value = 0
$wsAll.Range("P77").Value = 0.9
$wsAll.Range("Q77").Value = "Downward momentum: -1.980% over 10 samples"

# --- "All Trades" row 78: trade #77, HighProbConvergence, UP, OPEN ---
$wsAll.Range("A78").Value = 77
Set-TextCell $wsAll 78 2 "2026-02-18"
Set-TextCell $wsAll 78 3 "00:12:16"
$wsAll.Range("D78").Value = "HighProbConvergence"
$wsAll.Range("E78").Value = "UP"
$wsAll.Range("F78").Value = 0.35
Set-EmptyTextCell $wsAll 78 7               # G78 Exit Price (unset)
$wsAll.Range("H78").Value = "OPEN"
$wsAll.Range("I78").Value = 0
$wsAll.Range("J78").Value = 0
$wsAll.Range("K78").Value = 100
Set-EmptyTextCell $wsAll 78 12              # L78 Exit Reason (unset)
$wsAll.Range("M78").Value = 0
$wsAll.Range("N78").Value = 0
$wsAll.Range("O78").Value = 0
$wsAll.Range("P78").Value = 0.95
$wsAll.Range("Q78").Value = "Mean reversion UP: price 1.88% below mean (z=-4.36)"

# --- "momentum" sheet row 11: trade #76, DOWN, OPEN ---
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Range("A11").Value = 76
Set-TextCell $wsMomentum 11 2 "2026-02-18"
Set-TextCell $wsMomentum 11 3 "00:12:16"
$wsMomentum.Range("D11").Value = "momentum"
$wsMomentum.Range("E11").Value = "DOWN"
$wsMomentum.Range("F11").Value = 0.65
Set-EmptyTextCell $wsMomentum 11 7           # G11 Exit Price (unset)
$wsMomentum.Range("H11").Value = "OPEN"
$wsMomentum.Range("I11").Value = 0
$wsMomentum.Range("J11").Value = 0
$wsMomentum.Range("K11").Value = 100
$wsMomentum.Range("L11").Value = 0
$wsMomentum.Range("M11").Value = 0
$wsMomentum.Range("N11").Value = 0.9
$wsMomentum.Range("O11").Value = "Downward momentum: -1.980% over 10 samples"
Set-EmptyTextCell $wsMomentum 11 16          # P11 Exit Reason (unset)
$wsMomentum.Range("Q11").Value = 0

# --- "HighProbConvergence" sheet row 5: trade #77, UP, OPEN ---
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")
$wsHPC.Range("A5").Value = 77
Set-TextCell $wsHPC 5 2 "2026-02-18"
Set-TextCell $wsHPC 5 3 "00:12:16"
$wsHPC.Range("D5").Value = "HighProbConvergence"
$wsHPC.Range("E5").Value = "UP"
$wsHPC.Range("F5").Value = 0.35
Set-EmptyTextCell $wsHPC 5 7                 # G5 Exit Price (unset)
$wsHPC.Range("H5").Value = "OPEN"
$wsHPC.Range("I5").Value = 0
$wsHPC.Range("J5").Value = 0
$wsHPC.Range("K5").Value = 100
$wsHPC.Range("L5").Value = 0
$wsHPC.Range("M5").Value = 0
$wsHPC.Range("N5").Value = 0.95
$wsHPC.Range("O5").Value = "Mean reversion UP: price 1.88% below mean (z=-4.36)"
Set-EmptyTextCell $wsHPC 5 16                # P5 Exit Reason (unset)
$wsHPC.Range("Q5").Value = 0

# --- "EMAArbitrage" sheet row 6: trade #75, UP, OPEN ---
$wsEMA = $wb.Worksheets.Item("EMAArbitrage")
$wsEMA.Range("A6").Value = 75
Set-TextCell $wsEMA 6 2 "2026-02-18"
Set-TextCell $wsEMA 6 3 "00:11:14"
$wsEMA.Range("D6").Value = "EMAArbitrage"
$wsEMA.Range("E6").Value = "UP"
$wsEMA.Range("F6").Value = 0.46
Set-EmptyTextCell $wsEMA 6 7                 # G6 Exit Price (unset)
$wsEMA.Range("H6").Value = "OPEN"
$wsEMA.Range("I6").Value = 0
$wsEMA.Range("J6").Value = 0
$wsEMA.Range("K6").Value = 100
$wsEMA.Range("L6").Value = 0
$wsEMA.Range("M6").Value = 0
$wsEMA.Range("N6").Value = 0.7169
$wsEMA.Range("O6").Value = "EMA:up, RSI:50.0, ROC:21.69% | 2/3 UP"
Set-EmptyTextCell $wsEMA 6 16                # P6 Exit Reason (unset)
$wsEMA.Range("Q6").Value = 0
